$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $val)
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '30.515.66'
Set-TextValue 'E2' '  +0.42%  '
Set-TextValue 'D3' '1.914.71'
Set-TextValue 'E4' '  -0.08%  '
Set-TextValue 'D5' '244.10'
Set-TextValue 'E5' '  +1.10%  '
Set-TextValue 'D6' '0.9999'
Set-TextValue 'D7' '0.4879'
Set-TextValue 'E7' '  +4.31%  '
Set-TextValue 'D8' '0.2896'
Set-TextValue 'E8' '  +2.19%  '
Set-TextValue 'E9' '  -3.48%  '
Set-TextValue 'B10' 'Solana'
Set-TextValue 'C10' 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue 'D10' '19.04'
Set-TextValue 'E10' '  +5.42%  '
Set-TextValue 'B11' 'Litecoin'
Set-TextValue 'C11' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D11' '107.29'
Set-TextValue 'E11' '  +1.72%  '
Set-TextValue 'D12' '1.924.20'
Set-TextValue 'E12' '  +0.57%  '
Set-TextValue 'E13' '  -0.22%  '
Set-TextValue 'D14' '5.270'
Set-TextValue 'E14' '  +2.34%  '
Set-TextValue 'D15' '0.6649'
Set-TextValue 'E15' '  +1.85%  '
Set-TextValue 'D16' '273.07'
Set-TextValue 'E16' '  -3.90%  '
Set-TextValue 'D17' '30.506.58'
Set-TextValue 'E17' '  +0.41%  '
Set-TextValue 'D18' '0.9997'
Set-TextValue 'E18' '  -0.15%  '
Set-TextValue 'D19' '0.000007538'
Set-TextValue 'E19' '  -0.68%  '
Set-TextValue 'B20' 'Avalanche'
Set-TextValue 'C20' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D20' '12.83'
Set-TextValue 'E20' '  -0.88%  '
Set-TextValue 'B21' 'WrappedliquidstakedEther2.0'
Set-TextValue 'C21' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D21' '2.167.95'
Set-TextValue 'E21' '  +0.58%  '
Set-TextValue 'D22' '5.497'
Set-TextValue 'E22' '  +5.77%  '
Set-TextValue 'D23' '1.000'
Set-TextValue 'E23' '  -0.02%  '
Set-TextValue 'D24' '6.387'
Set-TextValue 'E24' '  +3.81%  '
Set-TextValue 'D25' '9.405'
Set-TextValue 'E25' '  +2.25%  '
Set-TextValue 'D26' '163.76'
Set-TextValue 'D27' '20.04'
Set-TextValue 'E27' '  -3.76%  '
Set-TextValue 'D28' '2.103'
Set-TextValue 'E28' '  +3.79%  '
Set-TextValue 'D29' '0.1048'
Set-TextValue 'E29' '  -1.65%  '
Set-TextValue 'D30' '1.401'
Set-TextValue 'E30' '  +2.30%  '
Set-TextValue 'D31' '4.121'
Set-TextValue 'E31' '  +0.48%  '
Set-TextValue 'D32' '4.039'
Set-TextValue 'E32' '  +2.38%  '
Set-TextValue 'E33' '  -1.16%  '
Set-TextValue 'D34' '0.7271'
Set-TextValue 'E34' '  -1.04%  '
Set-TextValue 'D35' '1.134'
Set-TextValue 'E35' '  -0.61%  '
Set-TextValue 'D36' '1.000'
Set-TextValue 'E36' '  +0.01%  '
Set-TextValue 'D37' '2.721'
Set-TextValue 'E37' '  +0.26%  '
Set-TextValue 'D38' '0.02033'
Set-TextValue 'E38' '  +1.88%  '
Set-TextValue 'E39' '  +0.25%  '
Set-TextValue 'D40' '110.67'
Set-TextValue 'E40' '  +2.44%  '
Set-TextValue 'D41' '2.013'
Set-TextValue 'E41' '  -1.44%  '
Set-TextValue 'D42' '0.4422'
Set-TextValue 'E42' '  +5.85%  '
Set-TextValue 'D43' '0.8664'
Set-TextValue 'E43' '  -0.32%  '
Set-TextValue 'D44' '5.866'
Set-TextValue 'E44' '  +1.13%  '
Set-TextValue 'D45' '0.9997'
Set-TextValue 'E45' '  -0.12%  '
Set-TextValue 'D46' '67.82'
Set-TextValue 'E46' '  +1.30%  '
Set-TextValue 'D47' '7.329'
Set-TextValue 'E47' '  +3.27%  '
Set-TextValue 'D48' '9.297'
Set-TextValue 'E48' '  +1.90%  '
Set-TextValue 'E49' '  +3.91%  '
Set-TextValue 'D50' '47.40'
Set-TextValue 'E50' '  -8.90%  '
Set-TextValue 'B51' 'NEARProtocol'
Set-TextValue 'C51' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D51' '1.460'
Set-TextValue 'E51' '  +6.86%  '
